$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 948.5
$ws.Range("I18").Value = 948.5
$ws.Range("K18").Value = 948.5
$ws.Range("M18").Value = -664.5
$ws.Range("H19").Value = 35717010
$ws.Range("I19").Value = 3689.3333
$ws.Range("K19").Value = 3689.3333
$ws.Range("M19").Value = -3514.3333
$ws.Range("H40").Value = 4269.392
$ws.Range("I40").Value = 3846.077
$ws.Range("K40").Value = 3846.077
$ws.Range("M40").Value = -3671.077
$ws.Range("H96").Value = 1163.1177
$ws.Range("I96").Value = 1350.6666
$ws.Range("K96").Value = 4051.9998
$ws.Range("M96").Value = -2678.9998
$ws.Range("H137").Value = 24393640
$ws.Range("I137").Value = 90912150
$ws.Range("K137").Value = 272736450
$ws.Range("M137").Value = -272733900
$ws.Range("H138").Value = 3325.8525
$ws.Range("I138").Value = 2048.2222
$ws.Range("J138").Value = 3860.6743
$ws.Range("K138").Value = 6144.6666
$ws.Range("L138").Value = 11582.0229
$ws.Range("M138").Value = -1004.6666
$ws.Range("N138").Value = -21862.0229

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2448.1792
$ws.Range("I32").Value = 2292.754
$ws.Range("K32").Value = 2292.754
$ws.Range("M32").Value = -2005.754
$ws.Range("H45").Value = 200001460
$ws.Range("I45").Value = 200001460
$ws.Range("K45").Value = 200001460
$ws.Range("M45").Value = -200001083
$ws.Range("H74").Value = 10419684
$ws.Range("I74").Value = 12347734
$ws.Range("K74").Value = 12347734
$ws.Range("M74").Value = -12346860
$ws.Range("H77").Value = 10419684
$ws.Range("I77").Value = 12347734
$ws.Range("K77").Value = 61738670
$ws.Range("M77").Value = -61734302
$ws.Range("H128").Value = 69999.5
$ws.Range("J128").Value = 69999.5
$ws.Range("L128").Value = 69999.5
$ws.Range("N128").Value = -79959.5
$ws.Range("H132").Value = 4032.544
$ws.Range("I132").Value = 3349.111
$ws.Range("K132").Value = 10047.333
$ws.Range("M132").Value = -7517.332999999999
$ws.Range("H140").Value = 81022.125
$ws.Range("J140").Value = 81255.28999999999
$ws.Range("L140").Value = 81255.28999999999
$ws.Range("N140").Value = -91615.28999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6093.231
$ws.Range("J20").Value = 6701.091
$ws.Range("L20").Value = 6701.091
$ws.Range("N20").Value = -7195.091
$ws.Range("H86").Value = 3547.2
$ws.Range("I86").Value = 2718.3333
$ws.Range("K86").Value = 2718.3333
$ws.Range("M86").Value = -1595.3333
$ws.Range("H89").Value = 3547.2
$ws.Range("I89").Value = 2718.3333
$ws.Range("K89").Value = 13591.6665
$ws.Range("M89").Value = -7975.666499999999
$ws.Range("H94").Value = 1031.238
$ws.Range("I94").Value = 974.8125
$ws.Range("J94").Value = 1211.8
$ws.Range("K94").Value = 974.8125
$ws.Range("L94").Value = 1211.8
$ws.Range("M94").Value = -523.8125
$ws.Range("N94").Value = -2113.8
$ws.Range("H99").Value = 2115
$ws.Range("I99").Value = 1925.8889
$ws.Range("K99").Value = 1925.8889
$ws.Range("M99").Value = -427.8888999999999
$ws.Range("H134").Value = 2117.8333
$ws.Range("J134").Value = 4133
$ws.Range("L134").Value = 12399
$ws.Range("N134").Value = -17469

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37088.727
$ws.Range("I31").Value = 4037.647
$ws.Range("J31").Value = 72205.5
$ws.Range("K31").Value = 4037.647
$ws.Range("L31").Value = 72205.5
$ws.Range("M31").Value = -3742.647
$ws.Range("N31").Value = -72795.5
$ws.Range("H34").Value = 37088.727
$ws.Range("I34").Value = 4037.647
$ws.Range("J34").Value = 72205.5
$ws.Range("K34").Value = 4037.647
$ws.Range("L34").Value = 72205.5
$ws.Range("M34").Value = -3835.647
$ws.Range("N34").Value = -72609.5
$ws.Range("H86").Value = 3951.125
$ws.Range("I86").Value = 4101.5
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 4101.5
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -2978.5
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 3951.125
$ws.Range("I89").Value = 4101.5
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 20507.5
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -14891.5
$ws.Range("N89").Value = -28732
$ws.Range("H134").Value = 2596.919
$ws.Range("I134").Value = 1438.6
$ws.Range("K134").Value = 4315.799999999999
$ws.Range("M134").Value = -1780.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1977.1
$ws.Range("I29").Value = 579.4
$ws.Range("J29").Value = 3374.8
$ws.Range("K29").Value = 1738.2
$ws.Range("L29").Value = 10124.4
$ws.Range("M29").Value = -1461.2
$ws.Range("N29").Value = -10678.4
$ws.Range("H141").Value = 6974.36
$ws.Range("I141").Value = 4380.3687
$ws.Range("K141").Value = 13141.1061
$ws.Range("M141").Value = -7961.106100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 37808
$ws.Range("J123").Value = 37808
$ws.Range("L123").Value = 37808
$ws.Range("N123").Value = -42708
$ws.Range("H132").Value = 3199.25
$ws.Range("I132").Value = 2771.4807
$ws.Range("K132").Value = 8314.4421
$ws.Range("M132").Value = -5784.4421
$ws.Range("H136").Value = 55443.777
$ws.Range("J136").Value = 55443.777
$ws.Range("L136").Value = 166331.331
$ws.Range("N136").Value = -171431.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8897.791999999999
$ws.Range("I7").Value = 6401.294
$ws.Range("K7").Value = 6401.294
$ws.Range("M7").Value = -6289.294
$ws.Range("H40").Value = 8068.6206
$ws.Range("I40").Value = 7678.0356
$ws.Range("K40").Value = 7678.0356
$ws.Range("M40").Value = -7542.0356
$ws.Range("H46").Value = 3733.3333
$ws.Range("I46").Value = 3850
$ws.Range("J46").Value = 3500
$ws.Range("K46").Value = 3850
$ws.Range("L46").Value = 3500
$ws.Range("M46").Value = -3662
$ws.Range("N46").Value = -3876
$ws.Range("H61").Value = 4606.423
$ws.Range("I61").Value = 3250.7368
$ws.Range("K61").Value = 3250.7368
$ws.Range("M61").Value = -3048.7368
$ws.Range("H113").Value = 4606.423
$ws.Range("I113").Value = 3250.7368
$ws.Range("K113").Value = 3250.7368
$ws.Range("M113").Value = -1080.7368
$ws.Range("H126").Value = 8897.791999999999
$ws.Range("I126").Value = 6401.294
$ws.Range("K126").Value = 19203.882
$ws.Range("M126").Value = -16733.882
$ws.Range("H132").Value = 3810
$ws.Range("I132").Value = 3248.3696
$ws.Range("J132").Value = 7500.7144
$ws.Range("K132").Value = 9745.1088
$ws.Range("L132").Value = 22502.1432
$ws.Range("M132").Value = -7215.1088
$ws.Range("N132").Value = -27562.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H100").Value = 605
$ws.Range("I100").Value = 697.5
$ws.Range("K100").Value = 1395
$ws.Range("M100").Value = -854
$ws.Range("H116").Value = 78947
$ws.Range("J116").Value = 78947
$ws.Range("L116").Value = 78947
$ws.Range("N116").Value = -88125
$ws.Range("H132").Value = 2329.577
$ws.Range("I132").Value = 1336.5084
$ws.Range("K132").Value = 4009.5252
$ws.Range("M132").Value = -1479.5252
